$wb = $excel.ActiveWorkbook

# Rename the "Collection_QRS_EQ5D-5L" worksheet to "CRF_QRS_EQ5D-5L"
$ws = $wb.Worksheets.Item("Collection_QRS_EQ5D-5L")
$ws.Name = "CRF_QRS_EQ5D-5L"
